# Applies the "Updated cryptos list" data refresh to Sheet1 of the workbook.
# For each changed row, Coin/Link (B/C) are updated only when the coin in that
# ranking slot changed; Price (D) and Volume(1h) (E) are always refreshed to the
# new scraped values. D-column values that look like plain decimal numbers are
# forced to Text format first so Excel keeps them as strings (matching the
# original "number-looking" text cells) instead of auto-converting to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "47.365.48"
$ws.Range("E2").Value = "  +2.56%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "2.505.29"
$ws.Range("E3").Value = "  +2.23%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  +0.09%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.07"
$ws.Range("E5").Value = "  +1.17%  "

# Row 6: Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.75"
$ws.Range("E6").Value = "  +4.72%  "

# Row 7: XRP
$ws.Range("E7").Value = "  +1.54%  "

# Row 8: USDC
$ws.Range("E8").Value = "  -0.03%  "

# Row 9: Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.539"
$ws.Range("E9").Value = "  +0.20%  "

# Row 10: Avalanche
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.42"

# Row 11: Dogecoin
$ws.Range("E11").Value = "  +1.11%  "

# Row 12: TRON
$ws.Range("E12").Value = "  +1.06%  "

# Row 13: Chainlink
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.45"
$ws.Range("E13").Value = "  +1.04%  "

# Row 14: Polkadot
$ws.Range("E14").Value = "  +2.05%  "

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.896.62"
$ws.Range("E15").Value = "  +2.19%  "

# Row 16: WrappedEther
$ws.Range("D16").Value = "2.508.30"
$ws.Range("E16").Value = "  +2.88%  "

# Row 17: Polygon
$ws.Range("E17").Value = "  +1.93%  "

# Row 18: WrappedBTC
$ws.Range("D18").Value = "47.294.43"
$ws.Range("E18").Value = "  +2.71%  "

# Row 19: InternetComputer(DFINITY)
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.90"
$ws.Range("E19").Value = "  +2.53%  "

# Row 20: Uniswap
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.67"
$ws.Range("E20").Value = "  +4.19%  "

# Row 21: ShibaInu
$ws.Range("E21").Value = "  +0.90%  "

# Row 22: ImmutableX
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.71"
$ws.Range("E22").Value = "  +14.47%  "

# Row 23: Litecoin
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.52"
$ws.Range("E23").Value = "  -0.83%  "

# Row 24: BitcoinCash
$ws.Range("E24").Value = "  +0.32%  "

# Row 25: PancakeSwap
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.60"
$ws.Range("E25").Value = "  +3.97%  "

# Row 26: EthereumClassic
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.07"
$ws.Range("E26").Value = "  +0.53%  "

# Row 27: Dai
$ws.Range("E27").Value = "  -0.01%  "

# Row 28: Toncoin (was Cosmos)
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.30"
$ws.Range("E28").Value = "  +0.36%  "

# Row 29: Cosmos (was Toncoin)
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.07"
$ws.Range("E29").Value = "  +4.21%  "

# Row 30: InjectiveProtocol
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.61"
$ws.Range("E30").Value = "  +5.89%  "

# Row 31: Kaspa
$ws.Range("E31").Value = "  +8.55%  "

# Row 32: OKB
$ws.Range("E32").Value = "  +1.27%  "

# Row 33: Celestia
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.99"
$ws.Range("E33").Value = "  +0.39%  "

# Row 34: Filecoin
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.46"
$ws.Range("E34").Value = "  +2.20%  "

# Row 35: Hedera
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0795"
$ws.Range("E35").Value = "  +4.31%  "

# Row 36: FirstDigitalUSD
$ws.Range("E36").Value = "  +0.25%  "

# Row 37: ARBITRUM
$ws.Range("E37").Value = "  +5.37%  "

# Row 38: RenderToken
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.73"
$ws.Range("E38").Value = "  +4.92%  "

# Row 39: LidoDAOToken
$ws.Range("E39").Value = "  +1.89%  "

# Row 40: Stellar
$ws.Range("E40").Value = "  +1.42%  "

# Row 41: WEMIXToken (was Monero)
$ws.Range("B41").Value = "WEMIXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.24"
$ws.Range("E41").Value = "  +0.55%  "

# Row 42: Monero (was WEMIXToken)
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "121.49"
$ws.Range("E42").Value = "  -4.26%  "

# Row 43: EnergySwap
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.19"
$ws.Range("E43").Value = "  +1.72%  "

# Row 44: VeChain
$ws.Range("E44").Value = "  +2.44%  "

# Row 45: Maker
$ws.Range("D45").Value = "2.001.34"
$ws.Range("E45").Value = "  +1.89%  "

# Row 46: NEARProtocol
$ws.Range("E46").Value = "  +4.30%  "

# Row 47: ApeXProtocol
$ws.Range("E47").Value = "  -0.86%  "

# Row 48: FraxShare (was Stacks)
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.10"
$ws.Range("E48").Value = "  -0.06%  "

# Row 49: Stacks (was FraxShare)
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.78"
$ws.Range("E49").Value = "  -3.89%  "

# Row 50: THORChain
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.24"
$ws.Range("E50").Value = "  +4.51%  "

# Row 51: MultiversX
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "57.50"
$ws.Range("E51").Value = "  +5.44%  "
